$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove frozen panes (row 2/col B freeze split) before touching rows
$excel.ActiveWindow.Split = $false

# Unmerge the title row (A1:I1) so individual cells can hold values again
$ws.Range("A1:I1").UnMerge()

# Delete row 2 (old header row with URL/Description) -- its content is promoted to row 1
$ws.Rows("2:2").Delete()

# Delete rows 12-22 (previously rows 13-23), shrinking the sheet to 11 rows
$ws.Rows("12:22").Delete()

# --- Row 1: new header row ---
$ws.Range("A1").Value = "URL"
$ws.Range("B1").Value = "Description"

$headerRange = $ws.Range("A1:I1")
$headerRange.Font.Name = "Helvetica"
$headerRange.Font.Size = 10
$headerRange.Font.Bold = $true
$headerRange.Interior.Color = 12566717
$headerRange.Borders.LineStyle = 1
$headerRange.NumberFormat = "0"
$ws.Range("A1").NumberFormat = "General"
$ws.Range("B1").NumberFormat = "General"

$ws.Rows("1:1").RowHeight = 20.55

# --- Rows 2-11: data rows ---
$dataRange = $ws.Range("A2:I11")
$dataRange.Font.Name = "Helvetica"
$dataRange.Font.Size = 10
$dataRange.Font.Bold = $false
$dataRange.Borders.LineStyle = 1
$dataRange.NumberFormat = "0"

$colA = $ws.Range("A2:A11")
$colA.Font.Bold = $true
$colA.Interior.Color = 14408667

$ws.Rows("2:2").RowHeight = 20.55
$ws.Rows("3:11").RowHeight = 20.35

# --- Column widths ---
$ws.Range("A1:I11").Columns.ColumnWidth = 9.125

Write-Output "done"
